$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.620.39'
$ws.Range('E2').Value = '  -3.09%  '
$ws.Range('D3').Value = '2.086.52'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5164'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4390'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09220'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.79'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.174'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.44'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('D13').Value = '2.080.65'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.207'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.03%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.737'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.85'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001160'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.34'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +10.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06638'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.008'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.209'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = '29.706.91'
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.70'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  -3.89%  '
$ws.Range('D26').Value = '2.327.80'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.90'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.30'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.520'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.98'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.148'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1053'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.623'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.196'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.971'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.118'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.25'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06722'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2275'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.48'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6847'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.285'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6640'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.18'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.305'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.623'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.218'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.61'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('B50').Value = 'WEMIXTOKEN'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.169'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07093'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.70%  '
